$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the text value in A1 (shared string content)
$ws.Range("A1").Value = "C:\Users\martin.asenov\Desktop\BBJ\BBJ_Workpackage_In sequence.pdf_Bookmarked.pdf"

# Widen column A to fit the new, longer text
$ws.Columns.Item(1).ColumnWidth = 98.8
